$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append a new row (row 12) describing a "NiLatticeMagicAi" ROI entry,
# mirroring the layout of the existing ROI rows (Name, Y1, Y2, X1, X2,
# ImageSizeY, ImageSizeX, Angle, SubRoiCenterSize, SubRoiNRowColumn, SubRoiSeparation).

$ws.Cells.Item(12, 1).Value = "NiLatticeMagicAi"
$ws.Cells.Item(12, 2).Value = 1016
$ws.Cells.Item(12, 3).Value = 1534
$ws.Cells.Item(12, 4).Value = 1755
$ws.Cells.Item(12, 5).Value = 1811
$ws.Cells.Item(12, 6).Value = 2160
$ws.Cells.Item(12, 7).Value = 2560
$ws.Cells.Item(12, 8).Value = 349.10000000000002
$ws.Cells.Item(12, 9).Value = "[897 1572 100 50;995 1591 100 50;1142 1619 150 50]"
$ws.Cells.Item(12, 10).Value = "[1 1]"
$ws.Cells.Item(12, 11).Value = "[100 100]"
